$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to carry a "Definitions" block (rows 10-16) explaining the
# vaccination-efficacy terms and a few homework-style questions. That block
# is being removed from the sheet.
#
# Rows 10, 12, 14 and 15 each held a single, plainly-formatted (column-
# default-styled) text cell - clearing both their content and formatting
# drops them out of the saved file entirely.
#
# Row 11 (B11) and row 16 (A16, merged A16:F16) keep distinctive cell
# formatting that must survive, so only their text is removed.

$ws.Range("A10").Clear()
$ws.Range("B12").Clear()
$ws.Range("A14").Clear()
$ws.Range("A15").Clear()

$ws.Range("B11").ClearContents()
$ws.Range("A16").Value = ""

$ws.Range("D11").Select()
